$d = $word.ActiveDocument

# --- Remove the stray _GoBack bookmark from the "Seriousness vs severity"
#     paragraph; it will be re-added at the very end of the newly
#     inserted "Rich Boyce..." paragraph below.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Append a new bullet paragraph: "Rich Boyce suggest to add both. ..."
#     (same ListParagraph / ilvl 1 / numId 1 formatting as its neighbours)
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null
$richBoycePara = $d.Paragraphs.Last
$richBoyceXml = '<w:p ' + $wNs + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Rich Boyce suggest to add both. Which was done on Feb 09, 2017</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$richBoycePara.Range.InsertXML($richBoyceXml) | Out-Null

# --- Append a further bullet paragraph: "Maria is looking to add ..."
$richBoycePara2 = $d.Paragraphs.Last
$richBoycePara2.Range.InsertParagraphAfter() | Out-Null
$mariaPara = $d.Paragraphs.Last
$mariaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Maria is looking to add examples from the literature.</w:t></w:r>' +
    '</w:p>'
$mariaPara.Range.InsertXML($mariaXml) | Out-Null
